$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) labels ---
$ws.Range("A1").Value = "leafwidth"
$ws.Range("B1").Value = "leafcolor"
$ws.Range("C1").Value = "pltheight"
$ws.Range("D1").Value = "lfprickle"

# --- Clear column C for data rows 2-21 (the "pltheight" answers are removed) ---
$ws.Range("C2:C21").Value = ""

# --- Delete the now-unused rows 22-32 (dimension shrinks from K32 to K21) ---
$ws.Rows("22:32").Delete()

# --- Apply the new row height (22.5pt, custom) to rows 1-21 ---
$ws.Rows("1:21").RowHeight = 22.5

# --- Update header/footer text ---
$ws.PageSetup.LeftHeader = "sample questionnaire 2"
$ws.PageSetup.RightHeader = "06/01/2021"

# Re-assert fit-to-page settings (setting PageSetup header properties resets
# unspecified PageSetup fields, so make sure these remain as before)
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
